$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 481.78946
$ws.Range("I39").Value = 33.6
$ws.Range("J39").Value = 641.8570999999999
$ws.Range("K39").Value = 100.8
$ws.Range("L39").Value = 1925.5713
$ws.Range("M39").Value = 195.2
$ws.Range("N39").Value = -2517.5713
$ws.Range("H113").Value = 2514.2856
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 654
$ws.Range("N113").Value = -8508
$ws.Range("H132").Value = 1374.7693
$ws.Range("I132").Value = 1075.5652
$ws.Range("J132").Value = 3668.6667
$ws.Range("K132").Value = 3226.6956
$ws.Range("L132").Value = 11006.0001
$ws.Range("M132").Value = -696.6956
$ws.Range("N132").Value = -16066.0001
$ws.Range("H137").Value = 1612.0526
$ws.Range("I137").Value = 1317.9333
$ws.Range("K137").Value = 3953.7999
$ws.Range("M137").Value = -1403.7999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7805.278
$ws.Range("I61").Value = 8973
$ws.Range("J61").Value = 1966.6666
$ws.Range("K61").Value = 8973
$ws.Range("L61").Value = 1966.6666
$ws.Range("M61").Value = -8761
$ws.Range("N61").Value = -2390.6666
$ws.Range("H74").Value = 2554.8235
$ws.Range("I74").Value = 2569.3333
$ws.Range("J74").Value = 2520
$ws.Range("K74").Value = 2569.3333
$ws.Range("L74").Value = 2520
$ws.Range("M74").Value = -1695.3333
$ws.Range("N74").Value = -4268
$ws.Range("H77").Value = 2554.8235
$ws.Range("I77").Value = 2569.3333
$ws.Range("J77").Value = 2520
$ws.Range("K77").Value = 12846.6665
$ws.Range("L77").Value = 12600
$ws.Range("M77").Value = -8478.666499999999
$ws.Range("N77").Value = -21336
$ws.Range("H122").Value = 2853365.2
$ws.Range("I122").Value = 4276041.5
$ws.Range("K122").Value = 12828124.5
$ws.Range("M122").Value = -12825674.5
$ws.Range("H132").Value = 3968.2903
$ws.Range("I132").Value = 1955.8667
$ws.Range("J132").Value = 5854.9375
$ws.Range("K132").Value = 5867.6001
$ws.Range("L132").Value = 17564.8125
$ws.Range("M132").Value = -3337.6001
$ws.Range("N132").Value = -22624.8125
$ws.Range("H136").Value = 7805.278
$ws.Range("I136").Value = 8973
$ws.Range("J136").Value = 1966.6666
$ws.Range("K136").Value = 26919
$ws.Range("L136").Value = 5899.9998
$ws.Range("M136").Value = -24369
$ws.Range("N136").Value = -10999.9998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5152.5713
$ws.Range("I134").Value = 6033.6
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 18100.8
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -15565.8
$ws.Range("N134").Value = -13920
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 449.8421
$ws.Range("I22").Value = 382.14285
$ws.Range("J22").Value = 639.4
$ws.Range("K22").Value = 382.14285
$ws.Range("L22").Value = 639.4
$ws.Range("M22").Value = -32.14285000000001
$ws.Range("N22").Value = -1339.4
$ws.Range("H31").Value = 4954.65
$ws.Range("I31").Value = 1492.05
$ws.Range("J31").Value = 8417.25
$ws.Range("K31").Value = 1492.05
$ws.Range("L31").Value = 8417.25
$ws.Range("M31").Value = -1197.05
$ws.Range("N31").Value = -9007.25
$ws.Range("H34").Value = 4954.65
$ws.Range("I34").Value = 1492.05
$ws.Range("J34").Value = 8417.25
$ws.Range("K34").Value = 1492.05
$ws.Range("L34").Value = 8417.25
$ws.Range("M34").Value = -1290.05
$ws.Range("N34").Value = -8821.25
$ws.Range("H58").Value = 1672.2972
$ws.Range("I58").Value = 924.5
$ws.Range("J58").Value = 2242.0476
$ws.Range("K58").Value = 924.5
$ws.Range("L58").Value = 2242.0476
$ws.Range("M58").Value = -721.5
$ws.Range("N58").Value = -2648.0476
$ws.Range("H94").Value = 2900.5173
$ws.Range("J94").Value = 2487.238
$ws.Range("L94").Value = 2487.238
$ws.Range("N94").Value = -3389.238
$ws.Range("H132").Value = 1534.2745
$ws.Range("I132").Value = 1351.275
$ws.Range("J132").Value = 2199.7273
$ws.Range("K132").Value = 4053.825
$ws.Range("L132").Value = 6599.1819
$ws.Range("M132").Value = -1523.825
$ws.Range("N132").Value = -11659.1819
$ws.Range("H134").Value = 2958.0715
$ws.Range("I134").Value = 3423.762
$ws.Range("J134").Value = 1561
$ws.Range("K134").Value = 10271.286
$ws.Range("L134").Value = 4683
$ws.Range("M134").Value = -7736.286
$ws.Range("N134").Value = -9753
$ws.Range("H136").Value = 1672.2972
$ws.Range("I136").Value = 924.5
$ws.Range("J136").Value = 2242.0476
$ws.Range("K136").Value = 2773.5
$ws.Range("L136").Value = 6726.1428
$ws.Range("M136").Value = -223.5
$ws.Range("N136").Value = -11826.1428
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2675.125
$ws.Range("I80").Value = 2619
$ws.Range("J80").Value = 2768.6667
$ws.Range("K80").Value = 2619
$ws.Range("L80").Value = 2768.6667
$ws.Range("M80").Value = -1621
$ws.Range("N80").Value = -4764.6667
$ws.Range("H83").Value = 2675.125
$ws.Range("I83").Value = 2619
$ws.Range("J83").Value = 2768.6667
$ws.Range("K83").Value = 13095
$ws.Range("L83").Value = 13843.3335
$ws.Range("M83").Value = -8103
$ws.Range("N83").Value = -23827.3335
$ws.Range("H122").Value = 3023143
$ws.Range("I122").Value = 1853569.6
$ws.Range("J122").Value = 16668167
$ws.Range("K122").Value = 5560708.800000001
$ws.Range("L122").Value = 50004501
$ws.Range("M122").Value = -5558258.800000001
$ws.Range("N122").Value = -50009401
$ws.Range("H132").Value = 4186.5386
$ws.Range("I132").Value = 6356
$ws.Range("J132").Value = 2830.625
$ws.Range("K132").Value = 19068
$ws.Range("L132").Value = 8491.875
$ws.Range("M132").Value = -16538
$ws.Range("N132").Value = -13551.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2642.8125
$ws.Range("I61").Value = 2299.4443
$ws.Range("J61").Value = 3084.2856
$ws.Range("K61").Value = 2299.4443
$ws.Range("L61").Value = 3084.2856
$ws.Range("M61").Value = -2097.4443
$ws.Range("N61").Value = -3488.2856
$ws.Range("H113").Value = 2642.8125
$ws.Range("I113").Value = 2299.4443
$ws.Range("J113").Value = 3084.2856
$ws.Range("K113").Value = 2299.4443
$ws.Range("L113").Value = 3084.2856
$ws.Range("M113").Value = -129.4443000000001
$ws.Range("N113").Value = -7424.2856
$ws.Range("H122").Value = 8930696
$ws.Range("I122").Value = 11906678
$ws.Range("J122").Value = 2752.5
$ws.Range("K122").Value = 35720034
$ws.Range("L122").Value = 8257.5
$ws.Range("M122").Value = -35717584
$ws.Range("N122").Value = -13157.5
$ws.Range("H132").Value = 18340320
$ws.Range("I132").Value = 28654800
$ws.Range("J132").Value = 3469.111
$ws.Range("K132").Value = 85964400
$ws.Range("L132").Value = 10407.333
$ws.Range("M132").Value = -85961870
$ws.Range("N132").Value = -15467.333
$ws.Range("H136").Value = 5679.8667
$ws.Range("I136").Value = 6586.273
$ws.Range("J136").Value = 3187.25
$ws.Range("K136").Value = 19758.819
$ws.Range("L136").Value = 9561.75
$ws.Range("M136").Value = -17208.819
$ws.Range("N136").Value = -14661.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2565.5454
$ws.Range("I96").Value = 2200
$ws.Range("J96").Value = 2870.1667
$ws.Range("K96").Value = 2200
$ws.Range("L96").Value = 2870.1667
$ws.Range("M96").Value = -827
$ws.Range("N96").Value = -5616.1667
$ws.Range("H107").Value = 55555844
$ws.Range("I107").Value = 58823800
$ws.Range("K107").Value = 176471400
$ws.Range("M107").Value = -176469480
$ws.Range("H132").Value = 1448.6957
$ws.Range("I132").Value = 1094.742
$ws.Range("J132").Value = 2180.2
$ws.Range("K132").Value = 3284.226
$ws.Range("L132").Value = 6540.599999999999
$ws.Range("M132").Value = -754.2259999999997
$ws.Range("N132").Value = -11600.6
$ws.Range("H136").Value = 4655.5625
$ws.Range("I136").Value = 8950.666999999999
$ws.Range("J136").Value = 2078.5
$ws.Range("K136").Value = 26852.001
$ws.Range("L136").Value = 6235.5
$ws.Range("M136").Value = -24302.001
$ws.Range("N136").Value = -11335.5
